$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price column (column D) cells so numeric-looking
# strings (e.g. "175.90", "64.497.62") are preserved exactly, including
# trailing zeros / multi-dot formats, instead of being coerced to floats.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Assign the new values
$ws.Range("D2").Value = "64.497.62"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "3.352.77"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "557.97"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "175.90"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "3.343.10"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "53.91"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "9.09"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "3.884.46"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.120"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "18.16"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "3.347.01"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "64.912.26"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").Value = "0.989"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "452.21"
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").Value = "4.90"
$ws.Range("E23").Value = "  +9.06%  "
$ws.Range("D24").Value = "4.13"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "14.32"
$ws.Range("E25").Value = "  +7.54%  "
$ws.Range("D26").Value = "86.75"
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "10.80"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "8.74"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "30.96"
$ws.Range("E30").Value = "  +4.40%  "
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  +0.43%  "
$ws.Range("D32").Value = "11.49"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "572.28"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Value = "60.97"
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("D35").Value = "0.108"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "3.63"
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("D38").Value = "0.141"
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("D39").Value = "35.43"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0741"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.370"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").Value = "3.065.65"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("D43").Value = "2.80"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("E44").Value = "  +2.35%  "
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "3.18"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "2.44"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "0.998"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "140.35"
$ws.Range("E49").Value = "  +5.06%  "
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").Value = "8.20"
$ws.Range("E51").Value = "  -0.38%  "

# Restore the default "Normal" style on column D cells so the underlying
# style index matches the original (unstyled) cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
